$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.421.41"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "1.646.26"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'0.9997"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'300.05"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").Value = "'0.3803"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").Value = "'50.34"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "'0.3492"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "'22.05"
$ws.Range("D14").Value = "'6.322"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "'7.258"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "1.648.12"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'94.78"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").Value = "'0.06972"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'6.612"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").Value = "'17.41"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").Value = "'12.43"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").Value = "23.443.25"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'2.434"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").Value = "'2.983"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").Value = "'21.08"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "'150.21"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").Value = "'5.164"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "'131.77"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").Value = "1.826.97"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'6.869"
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("D33").Value = "'2.134"
$ws.Range("E33").Value = "  -5.27%  "
$ws.Range("D34").Value = "'11.20"
$ws.Range("E34").Value = "  -8.16%  "
$ws.Range("D35").Value = "'0.9922"
$ws.Range("E35").Value = "  -6.20%  "
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("D37").Value = "'0.08789"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2425"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.906"
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").Value = "'0.06835"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "'12.79"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").Value = "'0.6843"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").Value = "'15.51"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").Value = "'0.9985"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "'0.6366"
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("D47").Value = "'2.244"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").Value = "'3.917"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("D50").Value = "'126.97"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "'1.228"
$ws.Range("E51").Value = "  +2.57%  "
